$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '67.110.43'
$ws.Range('E2').Value = '  -3.45%  '

$ws.Range('D3').Value = '3.729.27'
$ws.Range('E3').Value = '  +0.30%  '

$ws.Range('E4').Value = '  -0.01%  '

Set-TextValue 'D5' '589.56'
$ws.Range('E5').Value = '  -3.77%  '

Set-TextValue 'D6' '171.70'
$ws.Range('E6').Value = '  -3.78%  '

$ws.Range('D7').Value = '3.728.66'
$ws.Range('E7').Value = '  +0.27%  '

$ws.Range('E8').Value = '  +0.02%  '

Set-TextValue 'D9' '0.518'
$ws.Range('E9').Value = '  -1.96%  '

$ws.Range('E10').Value = '  -4.74%  '

Set-TextValue 'D11' '6.26'
$ws.Range('E11').Value = '  -4.47%  '

Set-TextValue 'D12' '0.457'
$ws.Range('E12').Value = '  -4.49%  '

Set-TextValue 'D13' '37.51'
$ws.Range('E13').Value = '  -5.32%  '

Set-TextValue 'D14' '0.0000242'
$ws.Range('E14').Value = '  -4.56%  '

$ws.Range('D15').Value = '4.348.74'
$ws.Range('E15').Value = '  +0.27%  '

$ws.Range('D16').Value = '3.725.37'
$ws.Range('E16').Value = '  +0.34%  '

$ws.Range('D17').Value = '67.178.61'
$ws.Range('E17').Value = '  -3.42%  '

$ws.Range('E18').Value = '  -4.73%  '

Set-TextValue 'D19' '7.05'
$ws.Range('E19').Value = '  -5.68%  '

Set-TextValue 'D20' '15.96'
$ws.Range('E20').Value = '  -1.80%  '

Set-TextValue 'D21' '483.64'
$ws.Range('E21').Value = '  -3.44%  '

Set-TextValue 'D22' '9.02'
$ws.Range('E22').Value = '  -1.34%  '

Set-TextValue 'D23' '0.717'
$ws.Range('E23').Value = '  -0.02%  '

Set-TextValue 'D24' '83.44'
$ws.Range('E24').Value = '  -3.01%  '

Set-TextValue 'D25' '2.35'
$ws.Range('E25').Value = '  -9.11%  '

Set-TextValue 'D26' '0.0000135'
$ws.Range('E26').Value = '  +0.84%  '

Set-TextValue 'D27' '12.13'
$ws.Range('E27').Value = '  -5.81%  '

Set-TextValue 'D28' '10.09'
$ws.Range('E28').Value = '  -10.00%  '

$ws.Range('E29').Value = '  -0.09%  '

Set-TextValue 'D30' '2.89'
$ws.Range('E30').Value = '  -0.41%  '

Set-TextValue 'D31' '2.38'
$ws.Range('E31').Value = '  -2.99%  '

Set-TextValue 'D32' '31.85'
$ws.Range('E32').Value = '  +5.01%  '

Set-TextValue 'D33' '7.64'
$ws.Range('E33').Value = '  -4.23%  '

$ws.Range('E34').Value = '  -5.03%  '

Set-TextValue 'D35' '0.998'
$ws.Range('E35').Value = '  -0.04%  '

Set-TextValue 'D36' '0.997'
$ws.Range('E36').Value = '  -4.50%  '

$ws.Range('E37').Value = '  -2.23%  '

Set-TextValue 'D38' '5.67'
$ws.Range('E38').Value = '  -6.66%  '

Set-TextValue 'D39' '0.321'
$ws.Range('E39').Value = '  -7.32%  '

Set-TextValue 'D40' '446.21'
$ws.Range('E40').Value = '  +2.76%  '

Set-TextValue 'D41' '48.60'
$ws.Range('E41').Value = '  -2.10%  '

Set-TextValue 'D42' '1.96'
$ws.Range('E42').Value = '  -4.34%  '

Set-TextValue 'D43' '2.84'
$ws.Range('E43').Value = '  -5.71%  '

Set-TextValue 'D44' '8.17'
$ws.Range('E44').Value = '  -4.46%  '

Set-TextValue 'D45' '41.08'
$ws.Range('E45').Value = '  -9.62%  '

$ws.Range('D46').Value = '2.788.84'
$ws.Range('E46').Value = '  -5.40%  '

$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D47' '139.49'
$ws.Range('E47').Value = '  +0.17%  '

$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D48' '1.00'
$ws.Range('E48').Value = '  +0.03%  '

Set-TextValue 'D49' '0.0345'

Set-TextValue 'D50' '25.70'
$ws.Range('E50').Value = '  -4.80%  '

$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D51' '2.28'
$ws.Range('E51').Value = '  -7.15%  '
